# Applies the cryptos.xlsx data refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a "Price" (column D) cell while keeping it as plain
# text (matching the source data, which stores prices as inline strings, not
# numbers). A leading apostrophe forces Excel to treat the value as text; we
# then restore the default "Normal" style so no numeric/text format is recorded
# on the cell (keeping cell styling identical to the original workbook).
function Set-TextPrice($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($cellRef).Style = "Normal"
}

$ws.Range("D2").Value = '69.726.94'
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").Value = '3.707.88'
$ws.Range("E3").Value = '  +0.48%  '
Set-TextPrice "D4" '0.999'
$ws.Range("E4").Value = '  -0.04%  '
Set-TextPrice "D5" '677.97'
$ws.Range("E5").Value = '  -1.10%  '
Set-TextPrice "D6" '162.01'
$ws.Range("E6").Value = '  +1.10%  '
$ws.Range("E7").Value = '  -0.10%  '
Set-TextPrice "D8" '0.497'
$ws.Range("E8").Value = '  +0.59%  '
$ws.Range("E9").Value = '  +1.84%  '
$ws.Range("E10").Value = '  +0.78%  '
$ws.Range("E11").Value = '  +1.96%  '
$ws.Range("E12").Value = '  +0.91%  '
Set-TextPrice "D13" '32.87'
$ws.Range("E13").Value = '  +1.20%  '
$ws.Range("D14").Value = '3.715.07'
$ws.Range("E14").Value = '  +0.84%  '
$ws.Range("D15").Value = '69.716.43'
$ws.Range("E15").Value = '  +0.45%  '
$ws.Range("E16").Value = '  +2.01%  '
Set-TextPrice "D17" '16.10'
$ws.Range("E17").Value = '  +1.60%  '
Set-TextPrice "D18" '6.51'
$ws.Range("E18").Value = '  +1.45%  '
Set-TextPrice "D19" '474.44'
$ws.Range("E19").Value = '  +0.87%  '
Set-TextPrice "D20" '9.84'
$ws.Range("E20").Value = '  -1.64%  '
Set-TextPrice "D21" '0.654'
$ws.Range("E21").Value = '  +0.52%  '
Set-TextPrice "D22" '80.48'
$ws.Range("E22").Value = '  +0.88%  '
$ws.Range("D23").Value = '3.853.45'
$ws.Range("E23").Value = '  +0.45%  '
$ws.Range("E24").Value = '  +3.20%  '
$ws.Range("E25").Value = '  +0.02%  '
Set-TextPrice "D26" '11.03'
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("E27").Value = '  -1.55%  '
$ws.Range("E28").Value = '  -0.29%  '
Set-TextPrice "D29" '1.76'
$ws.Range("E29").Value = '  +0.61%  '
$ws.Range("E30").Value = '  +1.42%  '
Set-TextPrice "D31" '6.61'
$ws.Range("E31").Value = '  +0.46%  '
$ws.Range("B32").Value = 'Kaspa'
$ws.Range("C32").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextPrice "D32" '0.167'
$ws.Range("E32").Value = '  +4.29%  '
$ws.Range("B33").Value = 'Binance-PegBSC-USD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextPrice "D33" '1.00'
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextPrice "D34" '26.99'
$ws.Range("E34").Value = '  +0.09%  '
$ws.Range("D35").Value = '3.695.83'
$ws.Range("E35").Value = '  +0.88%  '
Set-TextPrice "D36" '8.52'
$ws.Range("E36").Value = '  +4.13%  '
$ws.Range("E37").Value = '  +1.23%  '
$ws.Range("E39").Value = '  +0.07%  '
$ws.Range("E40").Value = '  -0.01%  '
Set-TextPrice "D41" '0.0907'
$ws.Range("E41").Value = '  +0.79%  '
$ws.Range("E42").Value = '  +0.31%  '
Set-TextPrice "D43" '166.94'
$ws.Range("E43").Value = '  +0.69%  '
Set-TextPrice "D44" '47.00'
$ws.Range("E44").Value = '  -1.13%  '
Set-TextPrice "D45" '2.80'
$ws.Range("E45").Value = '  +2.99%  '
$ws.Range("B46").Value = 'FLOKI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
Set-TextPrice "D46" '0.000280'
$ws.Range("E46").Value = '  -1.05%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextPrice "D47" '28.26'
$ws.Range("E47").Value = '  +0.98%  '
$ws.Range("E48").Value = '  -1.35%  '
$ws.Range("E49").Value = '  -0.09%  '
$ws.Range("E50").Value = '  +1.59%  '
$ws.Range("E51").Value = '  +2.15%  '
